# Update the "Resumo de Inscricoes" data on the "Inscricoes" sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Row 16: Tecnico Subsequente em Logistica - Campus Ribeirao das Neves
$ws.Range("E16").Value = 202
$ws.Range("F16").Value = 48
$ws.Range("H16").Value = 48

# Row 18: Tecnico Subsequente em Seguranca do Trabalho - Campus Santa Luzia
$ws.Range("E18").Value = 46
